$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Append new row 54 to the Logs sheet with the new test-mail entry.
$row = 54
$logs.Cells.Item($row, 1).Value = "Stuur je me even de datasheet van VentiQ-250?"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #13: Stuur je me even de datasheet van VentiQ-250?"
$logs.Cells.Item($row, 4).Value = "Documentatie / Datasheets"
$logs.Cells.Item($row, 5).Value = "Bedankt, we hebben dit doorgestuurd naar documentatie@bedrijf.nl."
$logs.Cells.Item($row, 6).Value = "2025-08-05 19:49:27"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Ja"
$logs.Cells.Item($row, 9).Value = "Nee"
$logs.Cells.Item($row, 10).Value = "Nee"

# Update the Dashboard count for "Documentatie / Datasheets" (row 9) from 1 to 2.
$dashboard.Cells.Item(9, 2).Value = 2

# Extend the conditional-formatting ranges that covered rows 2-53 so they
# also cover the newly appended row 54.
$logs.Range("D2:D53").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D54"))
$logs.Range("G2:G53").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G54"))
$logs.Range("H2:H53").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H54"))
$logs.Range("I2:I53").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I54"))
$logs.Range("J2:J53").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J54"))
